$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Acierto, profit 0.48
$ws.Range("G6").Value = "Acierto"
$ws.Range("H6").Value = 0.48

# Row 10: Fallo, profit -1
$ws.Range("G10").Value = "Fallo"
$ws.Range("H10").Value = -1

# Row 12: Fallo, profit -1
$ws.Range("G12").Value = "Fallo"
$ws.Range("H12").Value = -1

# Row 14: Fallo, profit -1
$ws.Range("G14").Value = "Fallo"
$ws.Range("H14").Value = -1

# Row 15: Fallo, profit -1
$ws.Range("G15").Value = "Fallo"
$ws.Range("H15").Value = -1
